$p = $ppt.ActivePresentation

# --- Slide 18: "Third Party solutions" table -------------------------------
# The cell listing "fabrikam.locationfinder.wsp" (row 4, col 1 of "Table 4")
# gets its run re-touched (PowerPoint marks it dirty/err on re-proof).
$s18 = $p.Slides.Item(18)
$tblShape = $s18.Shapes.Item(2)
$tbl = $tblShape.Table
$cell = $tbl.Cell(4, 1)
$cellTr = $cell.Shape.TextFrame.TextRange
$cellTr.Text = "fabrikam.locationfinder.wsp"

# --- Slide 19: "Branding" overview ------------------------------------------
# First bullet "Current Implementation:" becomes
# "Requirements and Current Implementation:" (split into two runs).
$s19 = $p.Slides.Item(19)
$bodyShape = $s19.Shapes.Item(2)
$bodyTr = $bodyShape.TextFrame.TextRange
$lead = $bodyTr.Characters(1, 8)
$lead.Text = "Requirements and Current "
